$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '28.416.96'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.13%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.822.26'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.36%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.003'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.12%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '315.51'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.20%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.003'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.18%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5201'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +1.54%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3867'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -1.49%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.08093'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +5.53%  '
$ws.Range("B10").Value = 'Polygon'
$ws.Range("C10").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.118'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.64%  '
$ws.Range("B11").Value = 'OKB'
$ws.Range("C11").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '41.91'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.35%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '6.394'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +1.28%  '
$ws.Range("B13").Value = 'Solana'
$ws.Range("C13").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '20.96'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.62%  '
$ws.Range("B14").Value = 'BinanceUSD'
$ws.Range("C14").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.003'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.09%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.422'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.56%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.824.95'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.21%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '94.54'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.87%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.00001104'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.19%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06640'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -1.00%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.69'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.00%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.003'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.23%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.042'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.72%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '28.443.36'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.16%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.41'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +2.23%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.245'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.51%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '159.24'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.46%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '20.92'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.58%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.026.15'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.63%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.414'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.55%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '125.03'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.38%  '
$ws.Range("E31").Value = '  +1.68%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.084'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -2.93%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.683'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.16%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.677'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.36%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.07391'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +5.03%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '12.43'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +10.74%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.2201'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.73%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02344'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.86%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.161'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.13%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '8.734'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -2.55%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.6337'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.96%  '
$ws.Range("E42").Value = '  -0.04%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.381'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.67%  '
$ws.Range("B44").Value = 'Decentraland'
$ws.Range("C44").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.6143'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +4.10%  '
$ws.Range("B45").Value = 'EnergySwap'
$ws.Range("C45").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '13.40'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.35%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.790'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +2.02%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '127.24'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +2.00%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.986'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.33%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.202'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.26%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.06892'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.55%  '
$ws.Range("B51").Value = 'Aave'
$ws.Range("C51").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '73.73'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.85%  '
